$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34 - this shifts the existing rows 34-45
# down to 35-46, preserving all of their data and formatting.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with this week's data (same market /
# product / variety / quality / unit / origin as the prior entry, new
# date + volume + prices).
$ws.Cells.Item(34, 1).Value = 5
$ws.Cells.Item(34, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(34, 3).Value = "Maule"
$ws.Cells.Item(34, 4).Value = 44489
$ws.Cells.Item(34, 5).Value = 7
$ws.Cells.Item(34, 6).Value = 100112022
$ws.Cells.Item(34, 7).Value = "Arveja Verde"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 200
$ws.Cells.Item(34, 11).Value = 18000
$ws.Cells.Item(34, 12).Value = 18000
$ws.Cells.Item(34, 13).Value = 18000
$ws.Cells.Item(34, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(34, 15).Value = "Región del Maule"
$ws.Cells.Item(34, 16).Value = 720
$ws.Cells.Item(34, 17).Value = 25
$ws.Cells.Item(34, 18).Value = "Hortaliza"
